# Commit: "Changed trade_id logic for exit"
#
# The MPWizard trade log picked up a new trade (MP41 /
# NIFTY09NOV23P19250) that produced several zero-value "exit" rows before
# finally recording its real exit on 2023-11-05 14:21:00. The running
# cash-ledger sheet (DTD) gets one matching settlement line per new trade
# row. Re-select cell A1 on MPWizard and make it the active sheet again.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# MPWizard: append rows 12-16 for trade MP41
# ---------------------------------------------------------------------
$mp = $wb.Worksheets.Item("MPWizard")

$tradeId = "MP41"
$symbol  = "NIFTY09NOV23P19250"
$signal  = "Long"
$entrySerial = 45235.59791666667
$exitSerial  = 45235.59791666667

# Four identical zero-value interim rows (12-15), entry/exit stored as
# real date serials exactly like the rest of the sheet.
for ($r = 12; $r -le 15; $r++) {
    $mp.Cells.Item($r, 1).Value  = $tradeId
    $mp.Cells.Item($r, 2).Value  = $symbol
    $mp.Cells.Item($r, 3).Value  = $signal
    $mp.Cells.Item($r, 4).Value  = $entrySerial
    $mp.Cells.Item($r, 5).Value  = $exitSerial
    $mp.Cells.Item($r, 6).Value  = 0
    $mp.Cells.Item($r, 7).Value  = 0
    $mp.Cells.Item($r, 8).Value  = 0
    $mp.Cells.Item($r, 9).Value  = 0
    $mp.Cells.Item($r, 10).Value = 0
    $mp.Cells.Item($r, 11).Value = 50
    $mp.Cells.Item($r, 12).Value = 0
    $mp.Cells.Item($r, 13).Value = 35.4
    $mp.Cells.Item($r, 14).Value = -35.4
}

# Row 16: the trade finally exits - entry/exit stored as literal text
# timestamps (matches how this sheet originally recorded MP165/MP166).
$row16Date = "2023-11-05 14:21:00"
$mp.Cells.Item(16, 1).Value  = $tradeId
$mp.Cells.Item(16, 2).Value  = $symbol
$mp.Cells.Item(16, 3).Value  = $signal
$mp.Cells.Item(16, 4).Value  = $row16Date
$mp.Cells.Item(16, 5).Value  = $row16Date
$mp.Cells.Item(16, 6).Value  = 0
$mp.Cells.Item(16, 7).Value  = 0
$mp.Cells.Item(16, 8).Value  = 0
$mp.Cells.Item(16, 9).Value  = 0
$mp.Cells.Item(16, 10).Value = 0
$mp.Cells.Item(16, 11).Value = 50
$mp.Cells.Item(16, 12).Value = 0
$mp.Cells.Item(16, 13).Value = 35.4
$mp.Cells.Item(16, 14).Value = -35.4

# Re-select A1 on MPWizard and make it the active sheet.
$mp.Range("A1").Select()
$mp.Activate()

# ---------------------------------------------------------------------
# DTD: append the matching ledger rows 15-19 (one per MP41 trade row)
# ---------------------------------------------------------------------
$dtd = $wb.Worksheets.Item("DTD")

$dtdDate = "05-Nov-23"
$dtdDay  = "Sunday"
$balances = @("₹44,399.81", "₹44,364.41", "₹44,329.01", "₹44,293.61", "₹44,258.21")

for ($i = 0; $i -lt 5; $i++) {
    $r = 15 + $i
    $dtd.Cells.Item($r, 1).Value = 13 + ($i + 1)
    $dtd.Cells.Item($r, 2).Value = $dtdDate
    $dtd.Cells.Item($r, 3).Value = $dtdDay
    $dtd.Cells.Item($r, 4).Value = $tradeId
    $dtd.Cells.Item($r, 5).Value = "MPWizard"
    $dtd.Cells.Item($r, 6).Value = "-₹35.40"
    $dtd.Cells.Item($r, 7).Value = $balances[$i]
}
